$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.763.81'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.394.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.83%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '507.64'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.20'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +4.28%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.553'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.391.44'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0979'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.66%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +6.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.67'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.816.52'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.46%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '56.712.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.92%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.81'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.03%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.407.42'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.82%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.17'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.96%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '310.70'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.26'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.64'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.60'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.35'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.51'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0730'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.67'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.88%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.48%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -3.99%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.994'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.91'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.85'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.60'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.816'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.44'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '132.77'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.40'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.40%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.68%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0910'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '248.06'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0486'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0211'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.28'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.47%  '
